$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "27.398.38"

Set-TextCell $ws.Range("D3") "1.848.11"

Set-TextCell $ws.Range("D4") "1.003"
Set-TextCell $ws.Range("E4") "  -0.74%  "

Set-TextCell $ws.Range("D5") "321.05"
Set-TextCell $ws.Range("E5") "  -0.13%  "

Set-TextCell $ws.Range("D6") "1.003"
Set-TextCell $ws.Range("E6") "  -0.53%  "

Set-TextCell $ws.Range("D7") "0.4452"
Set-TextCell $ws.Range("E7") "  -6.47%  "

Set-TextCell $ws.Range("D8") "0.3828"
Set-TextCell $ws.Range("E8") "  -5.81%  "

Set-TextCell $ws.Range("D9") "48.99"
Set-TextCell $ws.Range("E9") "  -8.10%  "

Set-TextCell $ws.Range("D10") "0.07820"
Set-TextCell $ws.Range("E10") "  -7.51%  "

Set-TextCell $ws.Range("D11") "1.013"
Set-TextCell $ws.Range("E11") "  -4.32%  "

Set-TextCell $ws.Range("D12") "21.43"
Set-TextCell $ws.Range("E12") "  -3.40%  "

Set-TextCell $ws.Range("D13") "1.841.88"
Set-TextCell $ws.Range("E13") "  -6.88%  "

Set-TextCell $ws.Range("D14") "5.828"
Set-TextCell $ws.Range("E14") "  -5.59%  "

Set-TextCell $ws.Range("D15") "7.079"
Set-TextCell $ws.Range("E15") "  -7.11%  "

Set-TextCell $ws.Range("D16") "1.002"
Set-TextCell $ws.Range("E16") "  -0.91%  "

Set-TextCell $ws.Range("D17") "85.30"
Set-TextCell $ws.Range("E17") "  -4.62%  "

Set-TextCell $ws.Range("D18") "0.00001023"
Set-TextCell $ws.Range("E18") "  -4.68%  "

Set-TextCell $ws.Range("D19") "0.06463"
Set-TextCell $ws.Range("E19") "  -2.55%  "

Set-TextCell $ws.Range("D20") "16.97"
Set-TextCell $ws.Range("E20") "  -9.08%  "

Set-TextCell $ws.Range("E21") "  -0.59%  "

Set-TextCell $ws.Range("D22") "5.459"
Set-TextCell $ws.Range("E22") "  -6.26%  "

Set-TextCell $ws.Range("D23") "27.382.95"
Set-TextCell $ws.Range("E23") "  -3.91%  "

Set-TextCell $ws.Range("D24") "10.75"
Set-TextCell $ws.Range("E24") "  -7.18%  "

Set-TextCell $ws.Range("D25") "2.263"
Set-TextCell $ws.Range("E25") "  -1.13%  "

Set-TextCell $ws.Range("D26") "2.062.43"
Set-TextCell $ws.Range("E26") "  -6.72%  "

Set-TextCell $ws.Range("D27") "151.54"
Set-TextCell $ws.Range("E27") "  -2.18%  "

Set-TextCell $ws.Range("D28") "19.27"
Set-TextCell $ws.Range("E28") "  -4.63%  "

Set-TextCell $ws.Range("D29") "2.029"
Set-TextCell $ws.Range("E29") "  -6.10%  "

Set-TextCell $ws.Range("D30") "5.454"
Set-TextCell $ws.Range("E30") "  -8.41%  "

Set-TextCell $ws.Range("D31") "119.32"
Set-TextCell $ws.Range("E31") "  -3.46%  "

Set-TextCell $ws.Range("B32") "ARBITRUM"
Set-TextCell $ws.Range("C32") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws.Range("D32") "1.476"
Set-TextCell $ws.Range("E32") "  +1.94%  "

Set-TextCell $ws.Range("B33") "Stellar"
Set-TextCell $ws.Range("C33") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D33") "0.09284"
Set-TextCell $ws.Range("E33") "  -3.29%  "

Set-TextCell $ws.Range("D34") "0.9224"
Set-TextCell $ws.Range("E34") "  -6.13%  "

Set-TextCell $ws.Range("D35") "3.592"
Set-TextCell $ws.Range("E35") "  -1.84%  "

Set-TextCell $ws.Range("D36") "5.204"
Set-TextCell $ws.Range("E36") "  -6.87%  "

Set-TextCell $ws.Range("D37") "0.02211"
Set-TextCell $ws.Range("E37") "  -5.31%  "

Set-TextCell $ws.Range("B38") "Hedera"
Set-TextCell $ws.Range("C38") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws.Range("D38") "0.05928"
Set-TextCell $ws.Range("E38") "  -4.78%  "

Set-TextCell $ws.Range("B39") "TrustWalletToken"
Set-TextCell $ws.Range("C39") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Range("D39") "1.207"
Set-TextCell $ws.Range("E39") "  -4.03%  "

Set-TextCell $ws.Range("D40") "8.281"
Set-TextCell $ws.Range("E40") "  -5.95%  "

Set-TextCell $ws.Range("D41") "1.002"
Set-TextCell $ws.Range("E41") "  -0.61%  "

Set-TextCell $ws.Range("D42") "0.5879"
Set-TextCell $ws.Range("E42") "  -5.51%  "

Set-TextCell $ws.Range("D43") "0.1842"
Set-TextCell $ws.Range("E43") "  -4.14%  "

Set-TextCell $ws.Range("D44") "10.20"
Set-TextCell $ws.Range("E44") "  -8.48%  "

Set-TextCell $ws.Range("D45") "1.257"
Set-TextCell $ws.Range("E45") "  -5.64%  "

Set-TextCell $ws.Range("D46") "0.5630"
Set-TextCell $ws.Range("E46") "  -5.64%  "

Set-TextCell $ws.Range("D47") "12.15"
Set-TextCell $ws.Range("E47") "  -6.53%  "

Set-TextCell $ws.Range("D48") "3.351"
Set-TextCell $ws.Range("E48") "  -1.44%  "

Set-TextCell $ws.Range("D49") "1.913"
Set-TextCell $ws.Range("E49") "  -6.93%  "

Set-TextCell $ws.Range("D50") "0.06837"
Set-TextCell $ws.Range("E50") "  +0.10%  "

Set-TextCell $ws.Range("D51") "1.003"
Set-TextCell $ws.Range("E51") "  -10.94%  "
